$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Insert new sheet "2022-Q1" right after "2021-Q4" ---
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "2022-Q1"

# Copy header-row formatting (bold/centered/bordered) from "2021-Q4"
$ws1.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Copy the index-column formatting from "2021-Q4" onto A2:A13
$ws1.Range("A2").Copy()
$newSheet.Range("A2:A13").PasteSpecial(-4122)

# --- Header row text ---
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Force text columns B, D:G to store values as text (matches source data format,
# and keeps leading zeros in fund codes like 090010)
$newSheet.Range("B2:B13").NumberFormat = "@"
$newSheet.Range("D2:G13").NumberFormat = "@"

# --- Data rows ---
# row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "510880"
$newSheet.Range("C2").Value = "华泰柏瑞上证红利ETF"
$newSheet.Range("D2").Value = "181.00"
$newSheet.Range("E2").Value = "97.22"
$newSheet.Range("F2").Value = "2.60"
$newSheet.Range("G2").Value = "4.7060"
$newSheet.Range("H2").Value = 7

# row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "090010"
$newSheet.Range("C3").Value = "大成中证红利指数A"
$newSheet.Range("D3").Value = "34.51"
$newSheet.Range("E3").Value = "93.73"
$newSheet.Range("F3").Value = "1.42"
$newSheet.Range("G3").Value = "0.4900"
$newSheet.Range("H3").Value = 8

# row 4
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "515180"
$newSheet.Range("C4").Value = "易方达中证红利ETF"
$newSheet.Range("D4").Value = "16.55"
$newSheet.Range("E4").Value = "99.58"
$newSheet.Range("F4").Value = "1.51"
$newSheet.Range("G4").Value = "0.2499"
$newSheet.Range("H4").Value = 8

# row 5
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "515080"
$newSheet.Range("C5").Value = "招商中证红利ETF"
$newSheet.Range("D5").Value = "9.06"
$newSheet.Range("E5").Value = "99.25"
$newSheet.Range("F5").Value = "1.50"
$newSheet.Range("G5").Value = "0.1359"
$newSheet.Range("H5").Value = 8

# row 6
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "009140"
$newSheet.Range("C6").Value = "永赢竞争力精选混合"
$newSheet.Range("D6").Value = "1.51"
$newSheet.Range("E6").Value = "88.85"
$newSheet.Range("F6").Value = "4.23"
$newSheet.Range("G6").Value = "0.0639"
$newSheet.Range("H6").Value = 5

# row 7
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "007801"
$newSheet.Range("C7").Value = "大成中证红利指数C"
$newSheet.Range("D7").Value = "3.87"
$newSheet.Range("E7").Value = "93.73"
$newSheet.Range("F7").Value = "1.42"
$newSheet.Range("G7").Value = "0.0550"
$newSheet.Range("H7").Value = 8

# row 8
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "515890"
$newSheet.Range("C8").Value = "博时中证红利ETF"
$newSheet.Range("D8").Value = "2.59"
$newSheet.Range("E8").Value = "98.55"
$newSheet.Range("F8").Value = "1.49"
$newSheet.Range("G8").Value = "0.0386"
$newSheet.Range("H8").Value = 8

# row 9
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "161907"
$newSheet.Range("C9").Value = "万家中证红利指数(LOF)"
$newSheet.Range("D9").Value = "1.34"
$newSheet.Range("E9").Value = "94.87"
$newSheet.Range("F9").Value = "1.43"
$newSheet.Range("G9").Value = "0.0192"
$newSheet.Range("H9").Value = 9

# row 10
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "006729"
$newSheet.Range("C10").Value = "万家中证500指数增强A"
$newSheet.Range("D10").Value = "1.04"
$newSheet.Range("E10").Value = "93.64"
$newSheet.Range("F10").Value = "1.33"
$newSheet.Range("G10").Value = "0.0138"
$newSheet.Range("H10").Value = 3

# row 11
$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "006730"
$newSheet.Range("C11").Value = "万家中证500指数增强C"
$newSheet.Range("D11").Value = "0.61"
$newSheet.Range("E11").Value = "93.64"
$newSheet.Range("F11").Value = "1.33"
$newSheet.Range("G11").Value = "0.0081"
$newSheet.Range("H11").Value = 3

# row 12
$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "001273"
$newSheet.Range("C12").Value = "民生加银新动力灵活配置混合A"
$newSheet.Range("D12").Value = "0.04"
$newSheet.Range("E12").Value = "68.44"
$newSheet.Range("F12").Value = "1.76"
$newSheet.Range("G12").Value = "0.0007"
$newSheet.Range("H12").Value = 7

# row 13
$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = "001274"
$newSheet.Range("C13").Value = "民生加银新动力灵活配置混合D"
$newSheet.Range("D13").Value = "0.04"
$newSheet.Range("E13").Value = "68.44"
$newSheet.Range("F13").Value = "1.76"
$newSheet.Range("G13").Value = "0.0007"
$newSheet.Range("H13").Value = 7

# --- Update the "总计" (Total) summary sheet: prepend a 2022-Q1 row ---
$wsTotal = $wb.Worksheets.Item(3)
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").Style = "Normal"
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 12
$wsTotal.Range("D2").Value = 5.78
